$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Seed the formats for the three new rows by copying the existing
# date/start-time/end-time formatting (xlPasteFormats) down from row 4/5
# so the new cells reuse the existing style records instead of minting
# brand-new numFmt entries.
$ws.Range("A4:C5").Copy() | Out-Null
$ws.Range("A6:C7").PasteSpecial(-4122) | Out-Null
$ws.Range("A4:C4").Copy() | Out-Null
$ws.Range("A8:C8").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Row 6: 2/27/2025, 12:30 PM - 1:30 PM
$ws.Range("A6").Value = 45715
$ws.Range("B6").Value = 0.52083333333333337
$ws.Range("C6").Value = 0.5625
$ws.Range("D6").Value = "Filled out the overleaf doc a bit."

# Row 7: 3/12/2025, 11:00 AM - 12:30 PM
$ws.Range("A7").Value = 45728
$ws.Range("B7").Value = 0.45833333333333331
$ws.Range("C7").Value = 0.52083333333333337
$ws.Range("D7").Value = "Filled out the overleaf doc more, and worked on the progress report"

# Row 8: 3/12/2025, 6:00 AM - 6:30 AM
$ws.Range("A8").Value = 45728
$ws.Range("B8").Value = 0.25
$ws.Range("C8").Value = 0.27083333333333331
$ws.Range("D8").Value = "Added another related work to the overleaf document"

# Column B widens (best-fit) to fit the new start-time values entered above
$ws.Columns.Item(2).ColumnWidth = 12.5

# Update the active selection cell to match the author's final cursor position
$ws.Range("F9").Select()
